$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new data row at row 11 (pushes existing rows 11-31 down
#    to 12-32, shifting merged cells automatically).
# ------------------------------------------------------------------
$ws.Rows("11:11").Insert()

# ------------------------------------------------------------------
# 2. Copy the formatting (styles/borders/fill/font/number-format) of
#    the row above (row 10) into the freshly inserted row 11 so the
#    new line looks like every other product row.
# ------------------------------------------------------------------
$ws.Range("A10:Q10").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Re-create the merged cells for row 11 (same layout as every
#    other product row: A:B, C:G, H:K, L:M, N:O).
# ------------------------------------------------------------------
$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

# ------------------------------------------------------------------
# 4. Fill in the values for the new item: "FAWAR FRUIT 6 SACHETS".
# ------------------------------------------------------------------
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "FAWAR FRUIT 6 SACHETS"
$ws.Range("H11").Value = "2:4"
$ws.Range("L11").Value = "1"
$ws.Range("N11").Value = "24.00"
$ws.Range("P11").Value = "24.0000"
$ws.Range("Q11").Value = "1:0"

# ------------------------------------------------------------------
# 5. Every product row that used to sit at old-row R (11..29) is now
#    one row further down (12..30); bump its running counter in
#    column A by 1 so the numbering (5,6,7,...) stays contiguous.
# ------------------------------------------------------------------
for ($r = 12; $r -le 30; $r++) {
    $cell = $ws.Range("A$r")
    $cell.Value = $cell.Value + 1
}

# ------------------------------------------------------------------
# 6. Update the grand-total cell (old row 30 -> now row 31) to add
#    the new line's sale price (24.00).
# ------------------------------------------------------------------
$totalCell = $ws.Range("P31")
$totalCell.Value = $totalCell.Value + 24

# ------------------------------------------------------------------
# 7. Refresh the "printed at" timestamp in the footer (old row 31 ->
#    now row 32).
# ------------------------------------------------------------------
$ws.Range("A32").Value = "Saturday, 26 July, 2025 3:41 PM"
